$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.03019297354432
$ws.Cells.Item(2, 4).Value = 1.033745539711229
$ws.Cells.Item(2, 5).Value = 0.9926147277508489
$ws.Cells.Item(2, 6).Value = 1.028773387736611
$ws.Cells.Item(2, 9).Value = 1.035870631516463
$ws.Cells.Item(2, 10).Value = 1.035335573056189
$ws.Cells.Item(2, 11).Value = 1.036547004498317
$ws.Cells.Item(2, 12).Value = 0.9955398523336033
$ws.Cells.Item(2, 13).Value = 1.031589229717407
$ws.Cells.Item(2, 14).Value = 1.015750620056057

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.031224342826947
$ws.Cells.Item(3, 4).Value = 1.034513373094421
$ws.Cells.Item(3, 5).Value = 0.9936372048519304
$ws.Cells.Item(3, 6).Value = 1.030437785616165
$ws.Cells.Item(3, 9).Value = 1.036169142026881
$ws.Cells.Item(3, 10).Value = 1.036007874898217
$ws.Cells.Item(3, 11).Value = 1.037123991942281
$ws.Cells.Item(3, 12).Value = 0.9963617723202692
$ws.Cells.Item(3, 13).Value = 1.033059326586002
$ws.Cells.Item(3, 14).Value = 1.015975596372707

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.03189086076553
$ws.Cells.Item(4, 4).Value = 1.035009457587172
$ws.Cells.Item(4, 5).Value = 0.9942998659930995
$ws.Cells.Item(4, 6).Value = 1.031513817662645
$ws.Cells.Item(4, 9).Value = 1.03636048219687
$ws.Cells.Item(4, 10).Value = 1.036441490726724
$ws.Cells.Item(4, 11).Value = 1.037495929848946
$ws.Cells.Item(4, 12).Value = 0.9968940712668345
$ws.Cells.Item(4, 13).Value = 1.034009180647991
$ws.Cells.Item(4, 14).Value = 1.016120636485677

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.032170864091072
$ws.Cells.Item(5, 4).Value = 1.035217831552289
$ws.Cells.Item(5, 5).Value = 0.9945786998346017
$ws.Cells.Item(5, 6).Value = 1.031965963360658
$ws.Cells.Item(5, 9).Value = 1.036440487485011
$ws.Cells.Item(5, 10).Value = 1.036623447060057
$ws.Cells.Item(5, 11).Value = 1.037651955630296
$ws.Cells.Item(5, 12).Value = 0.997117960005301
$ws.Cells.Item(5, 13).Value = 1.034408171867434
$ws.Cells.Item(5, 14).Value = 1.01618148378798

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.032217866161842
$ws.Cells.Item(6, 4).Value = 1.035252807928331
$ws.Cells.Item(6, 5).Value = 0.9946255319796338
$ws.Cells.Item(6, 6).Value = 1.032041868050353
$ws.Cells.Item(6, 9).Value = 1.036453895286881
$ws.Cells.Item(6, 10).Value = 1.036653978680802
$ws.Cells.Item(6, 11).Value = 1.037678133317903
$ws.Cells.Item(6, 12).Value = 0.9971555583673453
$ws.Cells.Item(6, 13).Value = 1.034475145250952
$ws.Cells.Item(6, 14).Value = 1.016191692849788

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.031894602970337
$ws.Cells.Item(7, 4).Value = 1.035012242596425
$ws.Cells.Item(7, 5).Value = 0.9943035907982488
$ws.Cells.Item(7, 6).Value = 1.031519860102985
$ws.Cells.Item(7, 9).Value = 1.036361552936439
$ws.Cells.Item(7, 10).Value = 1.036443923353729
$ws.Cells.Item(7, 11).Value = 1.037498015995683
$ws.Cells.Item(7, 12).Value = 0.9968970624462087
$ws.Cells.Item(7, 13).Value = 1.034014513266112
$ws.Cells.Item(7, 14).Value = 1.016121450031016

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.030541706165952
$ws.Cells.Item(8, 4).Value = 1.034005189425485
$ws.Cells.Item(8, 5).Value = 0.9929600610674301
$ws.Cells.Item(8, 6).Value = 1.029336079209632
$ws.Cells.Item(8, 9).Value = 1.035971890972602
$ws.Cells.Item(8, 10).Value = 1.035563073091403
$ws.Cells.Item(8, 11).Value = 1.036742292991845
$ws.Cells.Item(8, 12).Value = 0.995817528259106
$ws.Cells.Item(8, 13).Value = 1.032086349297817
$ws.Cells.Item(8, 14).Value = 1.015826762806938

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.028151154290647
$ws.Cells.Item(9, 4).Value = 1.032224813075779
$ws.Cells.Item(9, 5).Value = 0.9906006454969559
$ws.Cells.Item(9, 6).Value = 1.025480384780415
$ws.Cells.Item(9, 9).Value = 1.035271317681589
$ws.Cells.Item(9, 10).Value = 1.034000051172794
$ws.Cells.Item(9, 11).Value = 1.035399755379935
$ws.Cells.Item(9, 12).Value = 0.9939188001724441
$ws.Cells.Item(9, 13).Value = 1.028677665143775
$ws.Cells.Item(9, 14).Value = 1.015303372448707

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.026552905933566
$ws.Cells.Item(10, 4).Value = 1.031033927091402
$ws.Cells.Item(10, 5).Value = 0.989033133672735
$ws.Cells.Item(10, 6).Value = 1.022904313680585
$ws.Cells.Item(10, 9).Value = 1.034794853115983
$ws.Cells.Item(10, 10).Value = 1.032950648107623
$ws.Cells.Item(10, 11).Value = 1.034497363199871
$ws.Cells.Item(10, 12).Value = 0.9926553831429383
$ws.Cells.Item(10, 13).Value = 1.026397363510594
$ws.Cells.Item(10, 14).Value = 1.014951651612078

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.025859739087017
$ws.Cells.Item(11, 4).Value = 1.030517305705898
$ws.Cells.Item(11, 5).Value = 0.988355674866747
$ws.Cells.Item(11, 6).Value = 1.021787403941967
$ws.Cells.Item(11, 9).Value = 1.034586295091216
$ws.Cells.Item(11, 10).Value = 1.032494471102116
$ws.Cells.Item(11, 11).Value = 1.034104854237934
$ws.Cells.Item(11, 12).Value = 0.9921088820399291
$ws.Cells.Item(11, 13).Value = 1.025408008195592
$ws.Cells.Item(11, 14).Value = 1.014798683387979

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.025602095949878
$ws.Cells.Item(12, 4).Value = 1.030325264044924
$ws.Cells.Item(12, 5).Value = 0.9881042295826724
$ws.Cells.Item(12, 6).Value = 1.021372305748718
$ws.Cells.Item(12, 9).Value = 1.034508488991242
$ws.Cells.Item(12, 10).Value = 1.032324757494965
$ws.Cells.Item(12, 11).Value = 1.033958791972111
$ws.Cells.Item(12, 12).Value = 0.9919059725120875
$ws.Cells.Item(12, 13).Value = 1.025040213082487
$ws.Cells.Item(12, 14).Value = 1.014741762801651

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.025657369017526
$ws.Cells.Item(13, 4).Value = 1.030366464236688
$ws.Cells.Item(13, 5).Value = 0.9881581567098651
$ws.Cells.Item(13, 6).Value = 1.021461356203338
$ws.Cells.Item(13, 9).Value = 1.03452519399115
$ws.Cells.Item(13, 10).Value = 1.03236117385119
$ws.Cells.Item(13, 11).Value = 1.033990134928134
$ws.Cells.Item(13, 12).Value = 0.9919494934313052
$ws.Cells.Item(13, 13).Value = 1.025119120300309
$ws.Cells.Item(13, 14).Value = 1.014753977057866

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.025838445710342
$ws.Cells.Item(14, 4).Value = 1.030501434457336
$ws.Cells.Item(14, 5).Value = 0.9883348863814464
$ws.Cells.Item(14, 6).Value = 1.02175309650736
$ws.Cells.Item(14, 9).Value = 1.034579870522881
$ws.Cells.Item(14, 10).Value = 1.032480448024454
$ws.Cells.Item(14, 11).Value = 1.034092786139911
$ws.Cells.Item(14, 12).Value = 0.9920921077337197
$ws.Cells.Item(14, 13).Value = 1.025377612398315
$ws.Cells.Item(14, 14).Value = 1.014793980386111

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.025949990417238
$ws.Cells.Item(15, 4).Value = 1.030584574749766
$ws.Cells.Item(15, 5).Value = 0.9884438009545853
$ws.Cells.Item(15, 6).Value = 1.021932816807768
$ws.Cells.Item(15, 9).Value = 1.034613513672997
$ws.Cells.Item(15, 10).Value = 1.032553901054713
$ws.Cells.Item(15, 11).Value = 1.034155997520498
$ws.Cells.Item(15, 12).Value = 0.9921799884222134
$ws.Cells.Item(15, 13).Value = 1.025536837304604
$ws.Cells.Item(15, 14).Value = 1.014818614305424

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.026598885423888
$ws.Cells.Item(16, 4).Value = 1.03106819322015
$ws.Cells.Item(16, 5).Value = 0.9890781214508737
$ws.Cells.Item(16, 6).Value = 1.022978407838022
$ws.Cells.Item(16, 9).Value = 1.034808647027964
$ws.Cells.Item(16, 10).Value = 1.032980885435185
$ws.Cells.Item(16, 11).Value = 1.034523375370046
$ws.Cells.Item(16, 12).Value = 0.9926916645766087
$ws.Cells.Item(16, 13).Value = 1.026462981535111
$ws.Cells.Item(16, 14).Value = 1.014961789425708

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.027005619818844
$ws.Cells.Item(17, 4).Value = 1.031371296434388
$ws.Cells.Item(17, 5).Value = 0.989476357848556
$ws.Cells.Item(17, 6).Value = 1.023633883345874
$ws.Cells.Item(17, 9).Value = 1.034930447102206
$ws.Cells.Item(17, 10).Value = 1.03324824391874
$ws.Cells.Item(17, 11).Value = 1.034753347636896
$ws.Cells.Item(17, 12).Value = 0.9930127773699352
$ws.Cells.Item(17, 13).Value = 1.027043393973066
$ws.Cells.Item(17, 14).Value = 1.015051419387944

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.027242753729312
$ws.Cells.Item(18, 4).Value = 1.031547998881024
$ws.Cells.Item(18, 5).Value = 0.9897087662937556
$ws.Cells.Item(18, 6).Value = 1.024016071491747
$ws.Cells.Item(18, 9).Value = 1.035001274410683
$ws.Cells.Item(18, 10).Value = 1.033404018126106
$ws.Cells.Item(18, 11).Value = 1.034887316159787
$ws.Cells.Item(18, 12).Value = 0.9932001317071769
$ws.Cells.Item(18, 13).Value = 1.027381749263275
$ws.Cells.Item(18, 14).Value = 1.015103634348491

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.027323592130589
$ws.Cells.Item(19, 4).Value = 1.031608234192282
$ws.Cells.Item(19, 5).Value = 0.9897880325774034
$ws.Cells.Item(19, 6).Value = 1.024146364384995
$ws.Cells.Item(19, 9).Value = 1.035025387980316
$ws.Cells.Item(19, 10).Value = 1.033457104079454
$ws.Cells.Item(19, 11).Value = 1.034932967103756
$ws.Cells.Item(19, 12).Value = 0.9932640239640975
$ws.Cells.Item(19, 13).Value = 1.027497087796886
$ws.Cells.Item(19, 14).Value = 1.015121427347968

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.02696199219303
$ws.Cells.Item(20, 4).Value = 1.031338785907709
$ws.Cells.Item(20, 5).Value = 0.9894336180360679
$ws.Cells.Item(20, 6).Value = 1.023563571538926
$ws.Cells.Item(20, 9).Value = 1.034917401519181
$ws.Cells.Item(20, 10).Value = 1.033219576633701
$ws.Cells.Item(20, 11).Value = 1.03472869142564
$ws.Cells.Item(20, 12).Value = 0.9929783193494215
$ws.Cells.Item(20, 13).Value = 1.026981140869681
$ws.Cells.Item(20, 14).Value = 1.015041809635633

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.025785127854531
$ws.Cells.Item(21, 4).Value = 1.030461693117489
$ws.Cells.Item(21, 5).Value = 0.9882828385668249
$ws.Cells.Item(21, 6).Value = 1.021667192636912
$ws.Cells.Item(21, 9).Value = 1.034563778995392
$ws.Cells.Item(21, 10).Value = 1.032445332199034
$ws.Cells.Item(21, 11).Value = 1.034062565297802
$ws.Cells.Item(21, 12).Value = 0.9920501090198102
$ws.Cells.Item(21, 13).Value = 1.025301501393028
$ws.Cells.Item(21, 14).Value = 1.014782203206352

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.025044201029885
$ws.Cells.Item(22, 4).Value = 1.02990938774961
$ws.Cells.Item(22, 5).Value = 0.9875604150241495
$ws.Cells.Item(22, 6).Value = 1.020473538866208
$ws.Cells.Item(22, 9).Value = 1.034339484324153
$ws.Cells.Item(22, 10).Value = 1.031956975759524
$ws.Cells.Item(22, 11).Value = 1.03364220007288
$ws.Cells.Item(22, 12).Value = 0.9914670000341481
$ws.Cells.Item(22, 13).Value = 1.024243679366145
$ws.Cells.Item(22, 14).Value = 1.014618391502507

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.025437074656841
$ws.Cells.Item(23, 4).Value = 1.03020225554458
$ws.Cells.Item(23, 5).Value = 0.9879432794643023
$ws.Cells.Item(23, 6).Value = 1.021106446190878
$ws.Cells.Item(23, 9).Value = 1.034458573110444
$ws.Cells.Item(23, 10).Value = 1.032216011107096
$ws.Cells.Item(23, 11).Value = 1.033865190600472
$ws.Cells.Item(23, 12).Value = 0.991776070289318
$ws.Cells.Item(23, 13).Value = 1.024804621117196
$ws.Cells.Item(23, 14).Value = 1.01470528701378

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.026981705970298
$ws.Cells.Item(24, 4).Value = 1.03135347630081
$ws.Cells.Item(24, 5).Value = 0.9894529299347244
$ws.Cells.Item(24, 6).Value = 1.023595342847134
$ws.Cells.Item(24, 9).Value = 1.034923296925616
$ws.Cells.Item(24, 10).Value = 1.033232530675125
$ws.Cells.Item(24, 11).Value = 1.034739833032299
$ws.Cells.Item(24, 12).Value = 0.9929938892766442
$ws.Cells.Item(24, 13).Value = 1.027009270950861
$ws.Cells.Item(24, 14).Value = 1.015046152068735

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.028769961354959
$ws.Cells.Item(25, 4).Value = 1.032685778364954
$ws.Cells.Item(25, 5).Value = 0.9912096547607049
$ws.Cells.Item(25, 6).Value = 1.026478125643688
$ws.Cells.Item(25, 9).Value = 1.035454088999704
$ws.Cells.Item(25, 10).Value = 1.034405425204348
$ws.Cells.Item(25, 11).Value = 1.035748127033177
$ws.Cells.Item(25, 12).Value = 0.9944092447426414
$ws.Cells.Item(25, 13).Value = 1.029560242027048
$ws.Cells.Item(25, 14).Value = 1.015439171773438
